$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131, shifting rows 131:234 down to 132:235
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with its data
$ws.Range("A131").Value = 4
$ws.Range("B131").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C131").Value = "Los Lagos"
$ws.Range("D131").Value = 44673
$ws.Range("E131").Value = 10
$ws.Range("F131").Value = "Fruta"
$ws.Range("G131").Value = 100108
$ws.Range("H131").Value = "Tropicales y subtropicales"
$ws.Range("I131").Value = 100108005
$ws.Range("J131").Value = "Piña"
$ws.Range("K131").Value = "Caramelo"
$ws.Range("L131").Value = "Tercera"
$ws.Range("M131").Value = 400
$ws.Range("N131").Value = 17000
$ws.Range("O131").Value = 18000
$ws.Range("P131").Value = 17500
$ws.Range("Q131").Value = "$/caja 16 unidades"
$ws.Range("R131").Value = "Ecuador"
$ws.Range("S131").Value = 1094
$ws.Range("T131").Value = 16
